$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 2 (7a2f0b06 entry) status changes to "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 2 (7a2f0b06 entry) gets Latest Target File / Latest Handback File / Latest Handback DateTime filled in
$wsZhCn.Range("I2").Value = "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md"
$wsZhCn.Range("J2").Value = "7a2f0b06-cb68-4a3a-9c58-6f449971e259.9d07e703626588052f69160e70a2e62890191b56.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-07 17:01:38"

# --- de-de sheet: row 2 (7a2f0b06 entry) gets Latest Target File / Latest Handback File / Latest Handback DateTime filled in
$wsDeDe.Range("I2").Value = "7a2f0b06-cb68-4a3a-9c58-6f449971e259.md"
$wsDeDe.Range("J2").Value = "7a2f0b06-cb68-4a3a-9c58-6f449971e259.9d07e703626588052f69160e70a2e62890191b56.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-07 17:01:46"
